# Added ERP Module Emp Details pages
#
# 1. Rename the "EmployeeDetails" sheet to "EmployeeInfo".
# 2. Update the employee id sample value in that sheet from "emp-001" to "EMP078".
# 3. Move the active selection on that sheet from J11 to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeDetails")

# Rename the sheet
$ws.Name = "EmployeeInfo"

# Update the sample employee id value (row 2, column A -> emp_ID header in row 1)
$ws.Range("A2").Value = "EMP078"

# Make this sheet active and move the selection to D11
$ws.Activate()
$ws.Range("D11").Select()
